$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "10u"
$ws.Range("B23").Value = "C22"
$ws.Range("C23").Value = "Capacitor_SMD:C_0603_1608Metric"
$ws.Range("D23").Value = "C19702"
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D61"))
$n = $wb.Names.Item(1)
$n.RefersTo = "=BOM!`$A`$1:`$D`$61"
[void]$ws.Range("C41").Select()
